$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Birth"

$ws.Range("E2").NumberFormat = "mm-dd-yy"
$ws.Range("E2").Copy()
$ws.Range("E3:E6").PasteSpecial(-4122)

$ws.Range("E2").Value = (Get-Date -Year 2019 -Month 1 -Day 2).Date
$ws.Range("E3").Value = (Get-Date -Year 2019 -Month 1 -Day 3).Date
$ws.Range("E4").Value = (Get-Date -Year 2019 -Month 1 -Day 4).Date
$ws.Range("E5").Value = (Get-Date -Year 2019 -Month 1 -Day 5).Date
$ws.Range("E6").Value = (Get-Date -Year 2019 -Month 1 -Day 6).Date

$ws.Range("G6").Select()
